$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row/column value updates (crypto price + 1h volume refresh)
$ws.Cells.Item(2, 4).Value2 = '59.169.23'  # D2
$ws.Cells.Item(2, 5).Value2 = '  +0.55%  '  # E2
$ws.Cells.Item(3, 4).Value2 = '2.527.40'  # D3
$ws.Cells.Item(3, 5).Value2 = '  +0.60%  '  # E3
$ws.Cells.Item(4, 5).Value2 = '  -0.02%  '  # E4
$ws.Cells.Item(5, 4).Value2 = '''540.26'  # D5
$ws.Cells.Item(5, 5).Value2 = '  +1.54%  '  # E5
$ws.Cells.Item(6, 4).Value2 = '''138.31'  # D6
$ws.Cells.Item(6, 5).Value2 = '  -0.14%  '  # E6
$ws.Cells.Item(7, 5).Value2 = '  +0.03%  '  # E7
$ws.Cells.Item(8, 5).Value2 = '  +1.09%  '  # E8
$ws.Cells.Item(9, 4).Value2 = '2.526.56'  # D9
$ws.Cells.Item(9, 5).Value2 = '  +0.51%  '  # E9
$ws.Cells.Item(10, 5).Value2 = '  +1.85%  '  # E10
$ws.Cells.Item(11, 5).Value2 = '  -0.60%  '  # E11
$ws.Cells.Item(12, 4).Value2 = '''5.37'  # D12
$ws.Cells.Item(12, 5).Value2 = '  -1.29%  '  # E12
$ws.Cells.Item(13, 5).Value2 = '  -1.65%  '  # E13
$ws.Cells.Item(14, 4).Value2 = '2.975.80'  # D14
$ws.Cells.Item(14, 5).Value2 = '  +0.70%  '  # E14
$ws.Cells.Item(15, 4).Value2 = '''23.24'  # D15
$ws.Cells.Item(15, 5).Value2 = '  +1.02%  '  # E15
$ws.Cells.Item(16, 4).Value2 = '59.151.88'  # D16
$ws.Cells.Item(16, 5).Value2 = '  +0.64%  '  # E16
$ws.Cells.Item(17, 5).Value2 = '  +0.24%  '  # E17
$ws.Cells.Item(18, 4).Value2 = '2.527.51'  # D18
$ws.Cells.Item(18, 5).Value2 = '  +0.74%  '  # E18
$ws.Cells.Item(19, 4).Value2 = '''11.13'  # D19
$ws.Cells.Item(19, 5).Value2 = '  +1.03%  '  # E19
$ws.Cells.Item(20, 4).Value2 = '''4.31'  # D20
$ws.Cells.Item(20, 5).Value2 = '  +1.07%  '  # E20
$ws.Cells.Item(21, 4).Value2 = '''326.43'  # D21
$ws.Cells.Item(21, 5).Value2 = '  +1.36%  '  # E21
$ws.Cells.Item(22, 5).Value2 = '  +0.11%  '  # E22
$ws.Cells.Item(23, 5).Value2 = '  +2.82%  '  # E23
$ws.Cells.Item(24, 4).Value2 = '''65.45'  # D24
$ws.Cells.Item(24, 5).Value2 = '  +5.45%  '  # E24
$ws.Cells.Item(25, 5).Value2 = '  +0.01%  '  # E25
$ws.Cells.Item(26, 5).Value2 = '  +0.72%  '  # E26
$ws.Cells.Item(27, 4).Value2 = '''1.00'  # D27
$ws.Cells.Item(27, 5).Value2 = '  +0.64%  '  # E27
$ws.Cells.Item(28, 4).Value2 = '''7.68'  # D28
$ws.Cells.Item(28, 5).Value2 = '  -0.94%  '  # E28
$ws.Cells.Item(29, 4).Value2 = '0.0₃0781'  # D29
$ws.Cells.Item(29, 5).Value2 = '  +1.64%  '  # E29
$ws.Cells.Item(30, 4).Value2 = '''6.74'  # D30
$ws.Cells.Item(30, 5).Value2 = '  +0.97%  '  # E30
$ws.Cells.Item(31, 5).Value2 = '  +0.54%  '  # E31
$ws.Cells.Item(32, 4).Value2 = '''167.42'  # D32
$ws.Cells.Item(32, 5).Value2 = '  +2.41%  '  # E32
$ws.Cells.Item(33, 5).Value2 = '  +7.03%  '  # E33
$ws.Cells.Item(34, 2).Value2 = 'ImmutableX'  # B34
$ws.Cells.Item(34, 3).Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'  # C34
$ws.Cells.Item(34, 4).Value2 = '''1.48'  # D34
$ws.Cells.Item(34, 5).Value2 = '  +3.89%  '  # E34
$ws.Cells.Item(35, 2).Value2 = 'USDe'  # B35
$ws.Cells.Item(35, 3).Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'  # C35
$ws.Cells.Item(35, 4).Value2 = '''0.999'  # D35
$ws.Cells.Item(35, 5).Value2 = '  +0.00%  '  # E35
$ws.Cells.Item(36, 5).Value2 = '  +0.45%  '  # E36
$ws.Cells.Item(37, 5).Value2 = '  -1.85%  '  # E37
$ws.Cells.Item(38, 5).Value2 = '  +0.11%  '  # E38
$ws.Cells.Item(39, 4).Value2 = '''36.84'  # D39
$ws.Cells.Item(39, 5).Value2 = '  +0.08%  '  # E39
$ws.Cells.Item(40, 4).Value2 = '''0.825'  # D40
$ws.Cells.Item(40, 5).Value2 = '  +2.88%  '  # E40
$ws.Cells.Item(41, 5).Value2 = '  +0.45%  '  # E41
$ws.Cells.Item(42, 4).Value2 = '''287.75'  # D42
$ws.Cells.Item(42, 5).Value2 = '  +3.35%  '  # E42
$ws.Cells.Item(43, 5).Value2 = '  +0.66%  '  # E43
$ws.Cells.Item(44, 4).Value2 = '''0.998'  # D44
$ws.Cells.Item(45, 4).Value2 = '''131.92'  # D45
$ws.Cells.Item(45, 5).Value2 = '  +8.48%  '  # E45
$ws.Cells.Item(46, 5).Value2 = '  +2.53%  '  # E46
$ws.Cells.Item(47, 5).Value2 = '  +0.05%  '  # E47
$ws.Cells.Item(48, 5).Value2 = '  +0.26%  '  # E48
$ws.Cells.Item(49, 4).Value2 = '''0.0512'  # D49
$ws.Cells.Item(49, 5).Value2 = '  +0.40%  '  # E49
$ws.Cells.Item(50, 5).Value2 = '  +0.13%  '  # E50
$ws.Cells.Item(51, 5).Value2 = '  -0.75%  '  # E51
